# Insert a new data row at row 402 (shifting existing rows 402:485 down to 403:486)
# and populate it with the new weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("402:402").Insert()

$ws.Range("A402").Value = 6
$ws.Range("B402").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C402").Value = "Metropolitana"
$ws.Range("D402").Value = 44711
$ws.Range("E402").Value = 13
$ws.Range("F402").Value = 100112052
$ws.Range("G402").Value = "Albahaca"
$ws.Range("H402").Value = "Sin especificar"
$ws.Range("I402").Value = "Primera"
$ws.Range("J402").Value = 60
$ws.Range("K402").Value = 4500
$ws.Range("L402").Value = 5000
$ws.Range("M402").Value = 4708
$ws.Range("N402").Value = "$/paquete"
$ws.Range("O402").Value = "Región de Arica y Parinacota"
$ws.Range("P402").Value = 4708
$ws.Range("Q402").Value = 1
$ws.Range("R402").Value = "Hortaliza"
